$d = $word.ActiveDocument

# The document contains 4 paragraphs of the form:
#   <run1>&lt;id&gt;</run1><run2>p162v_N</run2><run3>&lt;/id&gt;</run3>
# Each needs to be merged into a single run:
#   <run1>&lt;id&gt;p162v_N&lt;/id&gt;</run1>
# keeping run1's original formatting/xml:space intact (so we delete the
# text belonging to run2+run3 and re-insert it right after run1, instead
# of doing a blanket Find/Replace which would rebuild run1's <w:t> and
# drop its xml:space="preserve").

for ($n = 1; $n -le 4; $n++) {
    $needle = "<id>p162v_" + $n + "</id>"
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        continue
    }

    $pStart = $rng.Start
    $tagOpenLen = 4           # "<id>"
    $idLen = ("p162v_" + $n).Length
    $tagCloseLen = 5          # "</id>"

    $tailStart = $pStart + $tagOpenLen
    $tailEnd = $tailStart + $idLen + $tagCloseLen

    $tailRange = $d.Range($tailStart, $tailEnd)
    $tailText = $tailRange.Text
    $tailRange.Delete()

    $insPoint = $d.Range($tailStart, $tailStart)
    $insPoint.InsertAfter($tailText)
}
